# Weekly COVID dashboard refresh: new "last updated" timestamp, updated case
# counts, and three re-ranked countries (Burkina Faso, Nigeria, Consejo Danes
# para los Refugiados) that moved up the "Casos totales" ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 23:22"

# Country,Casos totales,Nuevos casos,Casos activos,Recuperados,Casos criticos,Muertes hoy,Muertes
# Rows are in final sort order (descending Casos totales) as they should appear
# starting at row 4, immediately below the header row (row 3).
$countryData = @(
    ,@("Estados Unidos", 393782, 26778, 21561, 359529, 9169, 1821, 12692)
    ,@("España", 140617, 3942, 43208, 83497, 7069, 571, 13912)
    ,@("Italia", 135586, 3039, 24392, 94067, 3792, 604, 17127)
    ,@("Francia", 109069, 11059, 19337, 79404, 7131, 1417, 10328)
    ,@("Alemania", 107458, 4083, 36081, 69394, 4895, 173, 1983)
    ,@("China", 81740, 32, 77167, 1242, 211, 0, 3331)
    ,@("Iran", 62589, 2089, 27039, 31678, 3987, 133, 3872)
    ,@("Reino Unido", 55242, 3634, 135, 48948, 1559, 786, 6159)
    ,@("Turquia", 34109, 3892, 1582, 31802, 1474, 76, 725)
    ,@("Suiza", 22253, 596, 8704, 12728, 391, 56, 821)
    ,@("Belgica", 22194, 1380, 4157, 16002, 1260, 403, 2035)
    ,@("Paises Bajos", 19580, 777, 250, 17229, 1424, 234, 2101)
    ,@("Canada", 17847, 1180, 3942, 13530, 426, 52, 375)
    ,@("Brasil", 13831, 1648, 127, 13023, 296, 117, 681)
    ,@("Austria", 12639, 342, 4046, 8350, 243, 23, 243)
    ,@("Portugal", 12442, 712, 184, 11913, 271, 34, 345)
    ,@("Corea del Sur", 10331, 47, 6694, 3445, 55, 6, 192)
    ,@("Israel", 9248, 344, 770, 8413, 149, 8, 65)
    ,@("Suecia", 7693, 487, 205, 6897, 640, 114, 591)
    ,@("Rusia", 7497, 1154, 494, 6945, 8, 11, 58)
    ,@("Australia", 5919, 24, 2547, 3324, 93, 3, 48)
    ,@("Noruega", 5907, 42, 32, 5786, 78, 13, 89)
    ,@("Irlanda", 5709, 345, 25, 5474, 165, 36, 210)
    ,@("India", 5311, 533, 421, 4740, 0, 14, 150)
    ,@("Chile", 5116, 301, 898, 4175, 337, 6, 43)
    ,@("Dinamarca", 5071, 390, 1491, 3377, 127, 16, 203)
    ,@("Chequia", 4944, 122, 172, 4684, 86, 10, 88)
    ,@("Polonia", 4848, 435, 191, 4528, 50, 22, 129)
    ,@("Rumania", 4417, 360, 460, 3760, 274, 21, 197)
    ,@("Pakistan", 4035, 269, 429, 3549, 28, 4, 57)
    ,@("Malasia", 3963, 170, 1321, 2579, 92, 1, 63)
    ,@("Japon", 3906, 0, 592, 3222, 79, 0, 92)
    ,@("Filipinas", 3764, 104, 84, 3503, 1, 14, 177)
    ,@("Ecuador", 3747, 0, 100, 3456, 156, 0, 191)
    ,@("Luxemburgo", 2970, 127, 500, 2426, 35, 3, 44)
    ,@("Peru", 2954, 393, 1301, 1546, 109, 15, 107)
    ,@("Arabia Saudita", 2795, 190, 615, 2139, 41, 3, 41)
    ,@("Indonesia", 2738, 247, 204, 2313, 0, 12, 221)
    ,@("Serbia", 2447, 247, 118, 2268, 109, 3, 61)
    ,@("Mexico", 2439, 296, 633, 1681, 89, 31, 125)
    ,@("Emiratos Arabes Unidos", 2359, 283, 186, 2161, 1, 1, 12)
    ,@("Finlandia", 2308, 132, 300, 1974, 81, 7, 34)
    ,@("Tailandia", 2258, 38, 888, 1343, 61, 1, 27)
    ,@("Panama", 2100, 0, 14, 2031, 88, 0, 55)
    ,@("Catar", 2057, 225, 150, 1901, 37, 2, 6)
    ,@("Republica Dominicana", 1956, 128, 36, 1822, 147, 12, 98)
    ,@("Grecia", 1832, 77, 269, 1482, 90, 2, 81)
    ,@("Colombia", 1780, 201, 100, 1630, 76, 4, 50)
    ,@("Sudafrica", 1749, 63, 95, 1641, 7, 1, 13)
    ,@("Argentina", 1628, 0, 338, 1234, 96, 3, 56)
    ,@("Islandia", 1586, 24, 559, 1021, 11, 0, 6)
    ,@("Singapur", 1481, 106, 377, 1098, 29, 0, 6)
    ,@("Argelia", 1468, 45, 113, 1162, 46, 20, 193)
    ,@("Ucrania", 1462, 143, 28, 1389, 16, 7, 45)
    ,@("Egipto", 1450, 128, 276, 1080, 0, 9, 94)
    ,@("Croacia", 1282, 60, 167, 1097, 35, 2, 18)
    ,@("Marruecos", 1184, 64, 93, 1001, 1, 10, 90)
    ,@("Nueva Zelanda", 1160, 54, 241, 918, 4, 0, 1)
    ,@("Estonia", 1149, 41, 69, 1059, 12, 2, 21)
    ,@("Irak", 1122, 91, 373, 684, 0, 1, 65)
    ,@("Eslovenia", 1059, 38, 102, 921, 30, 6, 36)
    ,@("Moldavia", 1056, 91, 40, 994, 80, 3, 22)
    ,@("Hong Kong", 936, 21, 236, 696, 12, 0, 4)
    ,@("Lituania", 880, 37, 8, 857, 11, 0, 15)
    ,@("Bielorrusia", 861, 161, 54, 794, 31, 0, 13)
    ,@("Armenia", 853, 20, 87, 758, 30, 0, 8)
    ,@("Hungria", 817, 73, 71, 699, 17, 9, 47)
    ,@("Barein", 811, 55, 458, 348, 4, 1, 5)
    ,@("Bosnia y Herzegovina", 764, 90, 68, 663, 4, 4, 33)
    ,@("Kuwait", 743, 78, 105, 637, 23, 0, 1)
    ,@("Azerbaiyan", 717, 76, 44, 665, 23, 1, 8)
    ,@("Crucero", 712, 0, 619, 82, 10, 0, 11)
    ,@("Kazajistan", 697, 35, 51, 640, 16, 0, 6)
    ,@("Camerun", 658, 0, 43, 606, 0, 0, 9)
    ,@("Tunez", 623, 27, 25, 575, 39, 1, 23)
    ,@("Republica de Macedonia", 599, 29, 30, 543, 15, 3, 26)
    ,@("Eslovaquia", 581, 47, 13, 566, 3, 0, 2)
    ,@("Bulgaria", 577, 28, 42, 512, 21, 1, 23)
    ,@("Letonia", 548, 6, 16, 530, 5, 1, 2)
    ,@("Libano", 548, 7, 62, 467, 29, 0, 19)
    ,@("Principado de Andorra", 545, 20, 39, 484, 17, 1, 22)
    ,@("Uzbekistan", 504, 47, 30, 472, 8, 0, 2)
    ,@("Republica de Chipre", 494, 29, 47, 438, 13, 0, 9)
    ,@("Costa Rica", 483, 16, 24, 457, 14, 0, 2)
    ,@("Afganistan", 423, 56, 18, 391, 0, 3, 14)
    ,@("Uruguay", 415, 0, 123, 286, 14, 0, 6)
    ,@("Cuba", 396, 46, 27, 358, 15, 2, 11)
    ,@("Burkina Faso", 384, 20, 127, 238, 0, 1, 19)
    ,@("Albania", 383, 6, 131, 230, 7, 1, 22)
    ,@("Taiwan", 376, 3, 61, 310, 0, 0, 5)
    ,@("Oman", 371, 40, 67, 302, 3, 0, 2)
    ,@("Jordania", 353, 4, 138, 209, 5, 0, 6)
    ,@("Reunion", 349, 0, 40, 309, 4, 0, 0)
    ,@("Costa de Marfil", 349, 26, 41, 305, 0, 0, 3)
    ,@("Honduras", 305, 7, 6, 277, 10, 0, 22)
    ,@("Malta", 293, 52, 5, 288, 4, 0, 0)
    ,@("Ghana", 287, 73, 31, 251, 2, 0, 5)
    ,@("Niger", 278, 25, 26, 241, 0, 1, 11)
    ,@("San Marino", 277, 0, 35, 210, 14, 0, 32)
    ,@("Mauricio", 268, 24, 8, 253, 3, 0, 7)
    ,@("Estado de Palestina", 261, 7, 42, 218, 0, 0, 1)
    ,@("Nigeria", 254, 16, 44, 204, 2, 1, 6)
    ,@("Vietnam", 249, 4, 123, 126, 8, 0, 0)
    ,@("Montenegro", 241, 8, 4, 235, 7, 0, 2)
    ,@("Senegal", 237, 11, 105, 130, 1, 0, 2)
    ,@("Kirguistan", 228, 12, 33, 191, 5, 0, 4)
    ,@("Georgia", 196, 8, 46, 147, 6, 1, 3)
    ,@("Bolivia", 194, 11, 2, 178, 3, 3, 14)
    ,@("Sri Lanka", 185, 7, 42, 137, 5, 1, 6)
    ,@("Islas Feroe", 184, 1, 120, 64, 2, 0, 0)
    ,@("Consejo Danes para los Refugiados", 180, 19, 9, 153, 0, 0, 18)
    ,@("Kenia", 172, 14, 7, 159, 2, 0, 6)
    ,@("Venezuela", 165, 0, 65, 93, 6, 0, 7)
    ,@("Mayotte", 164, 0, 15, 147, 3, 0, 2)
    ,@("Banglades", 164, 41, 33, 114, 1, 5, 17)
    ,@("Martinica", 152, 1, 50, 98, 20, 0, 4)
    ,@("Isla de Man", 150, 11, 73, 76, 6, 0, 1)
    ,@("Guinea", 144, 16, 5, 139, 0, 0, 0)
    ,@("Guadalupe", 139, 0, 31, 101, 14, 0, 7)
    ,@("Brunei", 135, 0, 85, 49, 3, 0, 1)
    ,@("Paraguay", 115, 2, 15, 95, 1, 0, 5)
    ,@("Camboya", 115, 1, 58, 57, 1, 0, 0)
    ,@("Gibraltar", 113, 4, 60, 53, 0, 0, 0)
    ,@("Trinidad yTobago", 107, 2, 1, 98, 0, 0, 8)
    ,@("Ruanda", 105, 0, 7, 98, 0, 0, 0)
    ,@("Republica de Yibuti", 90, 0, 9, 81, 0, 0, 0)
    ,@("Madagascar", 88, 6, 2, 86, 6, 0, 0)
    ,@("Monaco", 79, 2, 4, 74, 4, 0, 1)
    ,@("El Salvador", 78, 9, 5, 69, 4, 0, 4)
    ,@("Guatemala", 77, 7, 17, 57, 3, 0, 3)
    ,@("Liechtenstein", 77, 0, 55, 21, 0, 0, 1)
    ,@("Aruba", 74, 3, 14, 60, 0, 0, 0)
    ,@("Guayana Francesa", 72, 0, 34, 38, 1, 0, 0)
    ,@("Togo", 65, 7, 23, 39, 0, 0, 3)
    ,@("Barbados", 63, 3, 6, 54, 4, 1, 3)
    ,@("Jamaica", 59, 0, 8, 48, 0, 0, 3)
    ,@("Mali", 56, 9, 12, 39, 0, 0, 5)
    ,@("Uganda", 52, 0, 0, 52, 0, 0, 0)
    ,@("Etiopia", 52, 8, 4, 46, 1, 0, 2)
    ,@("Polinesia Francesa", 47, 5, 0, 47, 0, 0, 0)
    ,@("Islas Caimanes", 45, 6, 6, 38, 0, 0, 1)
    ,@("Congo", 45, 0, 2, 38, 0, 0, 5)
    ,@("Macao", 44, 0, 10, 34, 1, 0, 0)
    ,@("San Martin (Parte Holandesa)", 40, 3, 1, 33, 2, 0, 6)
    ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
    ,@("Zambia", 39, 0, 7, 31, 0, 0, 1)
    ,@("Bermudas", 39, 0, 17, 20, 0, 0, 2)
    ,@("Guinea-Bisau", 33, 15, 0, 33, 0, 0, 0)
    ,@("Bahamas", 33, 0, 5, 22, 1, 1, 6)
    ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
    ,@("San Martin (Parte Francesa)", 32, 0, 7, 23, 6, 0, 2)
    ,@("Eritrea", 31, 0, 0, 31, 0, 0, 0)
    ,@("Guyana", 31, 0, 8, 18, 8, 1, 5)
    ,@("Gabon", 30, 6, 1, 28, 0, 0, 1)
    ,@("Benin", 26, 0, 5, 20, 0, 0, 1)
    ,@("Haiti", 25, 1, 0, 24, 0, 0, 1)
    ,@("Tanzania", 24, 0, 5, 18, 0, 0, 1)
    ,@("Birmania", 22, 0, 0, 21, 0, 0, 1)
    ,@("Libia", 20, 1, 1, 18, 0, 0, 1)
    ,@("Siria", 19, 0, 3, 14, 0, 0, 2)
    ,@("Maldivas", 19, 0, 13, 6, 0, 0, 0)
    ,@("Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0)
    ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
    ,@("Angola", 17, 1, 2, 13, 0, 0, 2)
    ,@("Guinea Ecuatorial", 16, 0, 3, 13, 0, 0, 0)
    ,@("Namibia", 16, 0, 3, 13, 0, 0, 0)
    ,@("Antigua y Barbuda", 15, 0, 0, 15, 1, 0, 0)
    ,@("Fiyi", 15, 1, 0, 15, 0, 0, 0)
    ,@("Dominica", 15, 0, 1, 14, 0, 0, 0)
    ,@("Mongolia", 15, 0, 4, 11, 0, 0, 0)
    ,@("Laos", 14, 2, 0, 14, 0, 0, 0)
    ,@("Santa Lucia", 14, 0, 1, 13, 0, 0, 0)
    ,@("Sudan", 14, 2, 2, 10, 0, 0, 2)
    ,@("Liberia", 14, 0, 3, 8, 0, 0, 3)
    ,@("Curazao", 13, 0, 7, 5, 0, 0, 1)
    ,@("Granada", 12, 0, 0, 12, 2, 0, 0)
    ,@("Seychelles", 11, 0, 0, 11, 0, 0, 0)
    ,@("San Cristobal y Nieves", 11, 1, 0, 11, 0, 0, 0)
    ,@("Zimbabue", 11, 1, 0, 9, 0, 1, 2)
    ,@("Groenlandia", 11, 0, 10, 1, 0, 0, 0)
    ,@("Mozambique", 10, 0, 1, 9, 0, 0, 0)
    ,@("Surinam", 10, 0, 0, 9, 0, 0, 1)
    ,@("Republica del Chad", 10, 1, 2, 8, 0, 0, 0)
    ,@("Suazilandia", 10, 0, 4, 6, 0, 0, 0)
    ,@("Nepal", 9, 0, 1, 8, 0, 0, 0)
    ,@("Montserrat", 9, 0, 0, 7, 0, 0, 2)
    ,@("Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0)
    ,@("Islas Turcas y Caicos", 8, 0, 0, 7, 0, 0, 1)
    ,@("Somalia", 8, 1, 1, 7, 0, 0, 0)
    ,@("Malaui", 8, 3, 0, 7, 1, 1, 1)
    ,@("Santa Sede", 7, 0, 0, 7, 0, 0, 0)
    ,@("Belice", 7, 0, 0, 6, 1, 0, 1)
    ,@("San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0)
    ,@("Cabo Verde", 7, 0, 1, 5, 0, 0, 1)
    ,@("Sierra Leona", 6, 0, 0, 6, 0, 0, 0)
    ,@("Botsuana", 6, 0, 0, 5, 0, 0, 1)
    ,@("Nicaragua", 6, 0, 0, 5, 0, 0, 1)
    ,@("San Bartolome", 6, 0, 1, 5, 0, 0, 0)
    ,@("Mauritania", 6, 0, 2, 3, 0, 0, 1)
    ,@("Butan", 5, 0, 2, 3, 0, 0, 0)
    ,@("Sahara Occidental", 4, 0, 0, 4, 0, 0, 0)
    ,@("Santo Tome y Principe", 4, 0, 0, 4, 0, 0, 0)
    ,@("Gambia", 4, 0, 2, 1, 0, 0, 1)
    ,@("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0)
    ,@("Anguila", 3, 0, 0, 3, 0, 0, 0)
    ,@("Burundi", 3, 0, 0, 3, 0, 0, 0)
    ,@("Papua Nueva Guinea", 2, 0, 0, 2, 0, 0, 0)
    ,@("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0)
    ,@("Islas Malvinas", 2, 0, 0, 2, 0, 0, 0)
    ,@("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
    ,@("Sudan del Sur", 1, 0, 0, 1, 0, 0, 0)
    ,@("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
)

$row = 4
foreach ($country in $countryData) {
    $ws.Cells.Item($row, 1).Value = $country[0]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $country[$col]
    }
    $row = $row + 1
}

